# Update LR-pairs data: add "Resolving-Mac" sending cluster and refresh
# all TPM-derived metrics for the Crlf1 -> Cntfr pair (rows 2-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Crlf1"
$ws.Range("C2").Value = "Cntfr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3711636666666667
$ws.Range("H2").Value = 1.113491
$ws.Range("I2").Value = 0.0322568113697387
$ws.Range("J2").Value = 0.03225681136973869
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.03844333333333334
$ws.Range("N2").Value = 0.11533
$ws.Range("O2").Value = 0.006478890266597937
$ws.Range("P2").Value = 0.006478890266597936
$ws.Range("Q2").Value = 0.01426876855888889
$ws.Range("R2").Value = 0.12841891703
$ws.Range("S2").Value = 0.0002089883412148857
$ws.Range("T2").Value = 0.0002089883412148857

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Crlf1"
$ws.Range("C3").Value = "Cntfr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3711636666666667
$ws.Range("H3").Value = 1.113491
$ws.Range("I3").Value = 0.0322568113697387
$ws.Range("J3").Value = 0.03225681136973869
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.666771333333333
$ws.Range("N3").Value = 17.000314
$ws.Range("O3").Value = 0.95502617622222
$ws.Range("P3").Value = 0.9550261762222199
$ws.Range("Q3").Value = 2.103299626241555
$ws.Range("R3").Value = 18.929696636174
$ws.Range("S3").Value = 0.03080609921956298
$ws.Range("T3").Value = 0.03080609921956297

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Crlf1"
$ws.Range("C4").Value = "Cntfr"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3711636666666667
$ws.Range("H4").Value = 1.113491
$ws.Range("I4").Value = 0.0322568113697387
$ws.Range("J4").Value = 0.03225681136973869
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2284146666666667
$ws.Range("N4").Value = 0.685244
$ws.Range("O4").Value = 0.03849493351118214
$ws.Range("P4").Value = 0.03849493351118213
$ws.Range("Q4").Value = 0.08477922520044444
$ws.Range("R4").Value = 0.763013026804
$ws.Range("S4").Value = 0.001241723808960835
$ws.Range("T4").Value = 0.001241723808960835

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Crlf1"
$ws.Range("C5").Value = "Cntfr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.123096
$ws.Range("H5").Value = 18.369288
$ws.Range("I5").Value = 0.5321413985496108
$ws.Range("J5").Value = 0.5321413985496107
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.03844333333333334
$ws.Range("N5").Value = 0.11533
$ws.Range("O5").Value = 0.006478890266597937
$ws.Range("P5").Value = 0.006478890266597936
$ws.Range("Q5").Value = 0.23539222056
$ws.Range("R5").Value = 2.11852998504
$ws.Range("S5").Value = 0.003447685727516887
$ws.Range("T5").Value = 0.003447685727516886

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Crlf1"
$ws.Range("C6").Value = "Cntfr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.123096
$ws.Range("H6").Value = 18.369288
$ws.Range("I6").Value = 0.5321413985496108
$ws.Range("J6").Value = 0.5321413985496107
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.666771333333333
$ws.Range("N6").Value = 17.000314
$ws.Range("O6").Value = 0.95502617622222
$ws.Range("P6").Value = 0.9550261762222199
$ws.Range("Q6").Value = 34.698184884048
$ws.Range("R6").Value = 312.283663956432
$ws.Range("S6").Value = 0.5082089650663792
$ws.Range("T6").Value = 0.508208965066379

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Crlf1"
$ws.Range("C7").Value = "Cntfr"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.123096
$ws.Range("H7").Value = 18.369288
$ws.Range("I7").Value = 0.5321413985496108
$ws.Range("J7").Value = 0.5321413985496107
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2284146666666667
$ws.Range("N7").Value = 0.685244
$ws.Range("O7").Value = 0.03849493351118214
$ws.Range("P7").Value = 0.03849493351118213
$ws.Range("Q7").Value = 1.398604931808
$ws.Range("R7").Value = 12.587444386272
$ws.Range("S7").Value = 0.02048474775571475
$ws.Range("T7").Value = 0.02048474775571474

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Crlf1"
$ws.Range("C8").Value = "Cntfr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.943314666666667
$ws.Range("H8").Value = 14.829944
$ws.Range("I8").Value = 0.4296098542617661
$ws.Range("J8").Value = 0.4296098542617661
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.03844333333333334
$ws.Range("N8").Value = 0.11533
$ws.Range("O8").Value = 0.006478890266597937
$ws.Range("P8").Value = 0.006478890266597936
$ws.Range("Q8").Value = 0.1900374935022222
$ws.Range("R8").Value = 1.71033744152
$ws.Range("S8").Value = 0.002783395103211115
$ws.Range("T8").Value = 0.002783395103211114

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Crlf1"
$ws.Range("C9").Value = "Cntfr"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.943314666666667
$ws.Range("H9").Value = 14.829944
$ws.Range("I9").Value = 0.4296098542617661
$ws.Range("J9").Value = 0.4296098542617661
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.666771333333333
$ws.Range("N9").Value = 17.000314
$ws.Range("O9").Value = 0.95502617622222
$ws.Range("P9").Value = 0.9550261762222199
$ws.Range("Q9").Value = 28.01263384471289
$ws.Range("R9").Value = 252.113704602416
$ws.Range("S9").Value = 0.4102886563829997
$ws.Range("T9").Value = 0.4102886563829996

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Crlf1"
$ws.Range("C10").Value = "Cntfr"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.943314666666667
$ws.Range("H10").Value = 14.829944
$ws.Range("I10").Value = 0.4296098542617661
$ws.Range("J10").Value = 0.4296098542617661
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2284146666666667
$ws.Range("N10").Value = 0.685244
$ws.Range("O10").Value = 0.03849493351118214
$ws.Range("P10").Value = 0.03849493351118213
$ws.Range("Q10").Value = 1.129125571815111
$ws.Range("R10").Value = 10.162130146336
$ws.Range("S10").Value = 0.01653780277555534
$ws.Range("T10").Value = 0.01653780277555533

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Crlf1"
$ws.Range("C11").Value = "Cntfr"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.06894633333333333
$ws.Range("H11").Value = 0.206839
$ws.Range("I11").Value = 0.005991935818884376
$ws.Range("J11").Value = 0.005991935818884375
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.03844333333333334
$ws.Range("N11").Value = 0.11533
$ws.Range("O11").Value = 0.006478890266597937
$ws.Range("P11").Value = 0.006478890266597936
$ws.Range("Q11").Value = 0.002650526874444445
$ws.Range("R11").Value = 0.02385474187
$ws.Range("S11").Value = [double]"3.882109465504952e-05"
$ws.Range("T11").Value = [double]"3.88210946550495e-05"

# Row 12
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Crlf1"
$ws.Range("C12").Value = "Cntfr"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.06894633333333333
$ws.Range("H12").Value = 0.206839
$ws.Range("I12").Value = 0.005991935818884376
$ws.Range("J12").Value = 0.005991935818884375
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 5.666771333333333
$ws.Range("N12").Value = 17.000314
$ws.Range("O12").Value = 0.95502617622222
$ws.Range("P12").Value = 0.9550261762222199
$ws.Range("Q12").Value = 0.3907031052717778
$ws.Range("R12").Value = 3.516327947446
$ws.Range("S12").Value = 0.005722455553278103
$ws.Range("T12").Value = 0.0057224555532781

# Row 13
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Crlf1"
$ws.Range("C13").Value = "Cntfr"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.06894633333333333
$ws.Range("H13").Value = 0.206839
$ws.Range("I13").Value = 0.005991935818884376
$ws.Range("J13").Value = 0.005991935818884375
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.2284146666666667
$ws.Range("N13").Value = 0.685244
$ws.Range("O13").Value = 0.03849493351118214
$ws.Range("P13").Value = 0.03849493351118213
$ws.Range("Q13").Value = 0.01574835374622222
$ws.Range("R13").Value = 0.141735183716
$ws.Range("S13").Value = 0.0002306591709512248
$ws.Range("T13").Value = 0.0002306591709512247

